# Update extracted_event_info.xlsx to reflect the new event data
# (US date/time formats, trimmed titles/descriptions, extra end-time cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write a literal text value into a cell without letting Excel
# auto-convert date-like strings (e.g. "11/22/2025") into date serials.
# We temporarily force a text number format, assign the value, then
# reset the cell style back to "Normal" so no stray style/numFmt is
# left behind on the cell.
# ---------------------------------------------------------------------
function Set-TextValue {
    param($cell, [string]$value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Shared text blocks
# ---------------------------------------------------------------------
$walkTalkTitle = "Walk&Talk"
$walkTalkDesc = "Gemeinsam spazieren gehen, quatschen, neue Leute kennenlernen oder vertraulich mit studentischen Hosts über mentale Gesundheit reden. Ein Angebot von und für Studierende."
$openingDesc = "Die Eröffnung der Sonderausstellung „Kunst / Macht. Rubens’ Medici-Zyklus und der gedruckte Kanon“ findet in der Schlosskirche von Schloss Hohentübingen statt. Die Ausstellung wurde von Ariane Koller und Anna Pawlak in Zusammenarbeit mit dem SFB 1391 Andere Ästhetik und dem Museum der Universität Tübingen MUT konzipiert. Sie wird von einer gleichnamigen Publikation begleitet."
$exhibitionDesc = "Vom 21. November 2025 bis 18. Januar 2026 zeigt die Graphische Sammlung des Kunsthistorischen Instituts auf Schloss Hohentübingen erstmals das druckgraphische Galeriewerk La Gallerie du Palais du Luxembourg (1710) mit 27 Kupferstichen nach dem Medici-Zyklus von Peter Paul Rubens. Die Ausstellung macht sichtbar, wie diese eindrucksvollen Druckgraphiken Rubens’ monumentale Bildfolge im frühen 18. Jahrhundert europaweit verbreiteten, neu interpretierten und kunsthistorisch prägend machten. Sie zeigt zugleich die Bedeutung des Galeriewerks als Vermittlungsmedium und als eigenständiges Kunstobjekt, das die Rezeption und Deutung des Zyklus über Jahrhunderte beeinflusste."

# ---------------------------------------------------------------------
# Rows 2-6: the five "Walk&Talk" occurrences
# ---------------------------------------------------------------------
$walkTalkRows = @(
    @{ Row = 2; Date = "11/22/2025" },
    @{ Row = 3; Date = "11/28/2025" },
    @{ Row = 4; Date = "12/12/2025" },
    @{ Row = 5; Date = "12/13/2025" },
    @{ Row = 6; Date = "01/16/2026" }
)

foreach ($entry in $walkTalkRows) {
    $r = $entry.Row

    # A: Title - trimmed down
    $ws.Cells.Item($r, 1).Value = $walkTalkTitle

    # B / C: Start_Date / End_Date - now in US MM/DD/YYYY text form
    Set-TextValue $ws.Cells.Item($r, 2) $entry.Date
    Set-TextValue $ws.Cells.Item($r, 3) $entry.Date

    # D / E: Start_Time / End_Time - now in US hh:mm AM/PM form
    $ws.Cells.Item($r, 4).Value = "03:00 PM"
    $ws.Cells.Item($r, 5).Value = "04:30 PM"

    # F: Description - reworded / shortened
    $ws.Cells.Item($r, 6).Value = $walkTalkDesc
}

# ---------------------------------------------------------------------
# Row 7: exhibition opening event
# ---------------------------------------------------------------------
Set-TextValue $ws.Cells.Item(7, 2) "11/20/2025"
Set-TextValue $ws.Cells.Item(7, 3) "11/20/2025"
$ws.Cells.Item(7, 4).Value = "07:00 PM"
$ws.Cells.Item(7, 6).Value = $openingDesc

# ---------------------------------------------------------------------
# Row 8: exhibition run event
# ---------------------------------------------------------------------
Set-TextValue $ws.Cells.Item(8, 2) "11/21/2025"
Set-TextValue $ws.Cells.Item(8, 3) "01/18/2026"
$ws.Cells.Item(8, 4).Value = "10:00 AM"
$ws.Cells.Item(8, 5).Value = "05:00 PM"
$ws.Cells.Item(8, 6).Value = $exhibitionDesc

# New Registration_Needed cell (boolean FALSE) added for row 8
$ws.Cells.Item(8, 10).Value = $false
